# Insert a new data row at row 165 (shifting the existing rows 165:256 down
# to 166:257) and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(165).Insert()

$ws.Range("A165").Value2 = 3
$ws.Range("B165").Value2 = "Femacal de La Calera"
$ws.Range("C165").Value2 = "Coquimbo"
$ws.Range("D165").Value2 = 44529
$ws.Range("E165").Value2 = 5
$ws.Range("F165").Value2 = 100112031
$ws.Range("G165").Value2 = "Poroto verde"
$ws.Range("H165").Value2 = "Magnum"
$ws.Range("I165").Value2 = "Primera"
$ws.Range("J165").Value2 = 76
$ws.Range("K165").Value2 = 27000
$ws.Range("L165").Value2 = 28000
$ws.Range("M165").Value2 = 27500
$ws.Range("N165").Value2 = "$/malla 25 kilos"
$ws.Range("O165").Value2 = "Provincia de Limarí"
$ws.Range("P165").Value2 = 1100
$ws.Range("Q165").Value2 = 25
$ws.Range("R165").Value2 = "Hortaliza"

# Keep the date style (same as the rest of column D) on the new cell.
$ws.Range("D165").NumberFormat = $ws.Range("D166").NumberFormat
